$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C5-ImproveRSSI row (row 10): update Pages, Finish, and Comment
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = "YES"
$ws.Range("E10").Value = "CHECKING LATER"

# Move the active selection to the SUM cell
$ws.Range("C11").Select()
